$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values
# G2 / G3 hold numeric-looking text; force Text format so the stored
# value stays a string (matching the original "3876.1"/"15.7" text cells)
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "3876.0"

$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "15.5"

$ws.Range("D4").Value = "LTH0330"
$ws.Range("I4").Value = "SCECO+STB"
$ws.Range("J4").Value = "Good"

# Remove row 5 entirely (shift cells up)
$ws.Rows.Item(5).Delete()
